# "add default penalties to amortising EAD"
#
# Inserts two new assumption columns (ead_default_penalty_pct,
# ead_default_penalty_amt) into the ASSUMPTIONS sheet, right before the
# existing eir_base_rate column, and populates row 3 (Segment 2) with a
# 2% default penalty and a 1 (flat amount) default penalty flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASSUMPTIONS")
$tm = $wb.Worksheets.Item("TRANSITION_MATRIX")

# --- 1. Insert two new columns where the old "eir_base_rate" (AE) column
#        used to be. This shifts eir_base_rate from AE -> AG and copies
#        the left neighbour's (AD) formatting into the new AE/AF cells,
#        matching how Excel seeds formatting for inserted columns.
$ws.Range("AE1:AF1").EntireColumn.Insert()

# --- 2. Header row: new field names (appended as new shared strings).
$ws.Range("AE1").Value = "ead_default_penalty_pct"
$ws.Range("AF1").Value = "ead_default_penalty_amt"

# Give the two new headers the same look (bold/white on blue) as every
# other header cell.
$ws.Range("AD1").Copy()
$ws.Range("AE1:AF1").PasteSpecial(-4122)
$ws.Range("AE1").Value = "ead_default_penalty_pct"
$ws.Range("AF1").Value = "ead_default_penalty_amt"

# --- 3. Data rows: copy the percentage-style formatting already used by
#        the neighbouring ead_prepayment_pct column (AD, style index 3)
#        down both new columns for every data row.
$ws.Range("AD2:AD6").Copy()
$ws.Range("AE2:AE6").PasteSpecial(-4122)
$ws.Range("AD2:AD6").Copy()
$ws.Range("AF2:AF6").PasteSpecial(-4122)

# --- 4. Values: only "Segment 2" (row 3) gets a default penalty assumption.
$ws.Range("AE3").Value = 0.02
$ws.Range("AF3").Value = 1

# AF3 (the penalty amount) is displayed with the accounting "Comma" style.
$ws.Range("AF3").Style = "Comma"

# --- 5. Column widths: widen the two new columns to fit their longer
#        header text (stored widths run ~0.83 wider than the character
#        count supplied to ColumnWidth).
$ws.Columns("AE").ColumnWidth = 22.74
$ws.Columns("AF").ColumnWidth = 23.45

# --- 6. View state: ASSUMPTIONS becomes the active/selected tab (moving
#        away from TRANSITION_MATRIX), with a new selected cell.
$ws.Activate() | Out-Null
$ws.Range("AF17").Select() | Out-Null
